$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.245.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.356.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.22%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.651"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.45"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +13.94%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.499"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +13.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0979"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "27.37"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.708.47"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.39%  "
$ws.Range("E13").Value = "  +2.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.43%  "
$ws.Range("E15").Value = "  +4.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.864"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.349.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.236.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "75.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.14%  "
$ws.Range("E21").Value = "  +4.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "250.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.14%  "
$ws.Range("E23").Value = "  +5.11%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("E27").Value = "  +3.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  +10.09%  "
$ws.Range("E31").Value = "  +3.32%  "
$ws.Range("E32").Value = "  +3.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("E34").Value = "  +2.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.27%  "
$ws.Range("E36").Value = "  +4.49%  "
$ws.Range("E37").Value = "  +7.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.97%  "
$ws.Range("E39").Value = "  +2.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +14.00%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.51%  "
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.08%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.50%  "
$ws.Range("E46").Value = "  +1.89%  "
$ws.Range("E47").Value = "  +2.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.443.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.583.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.62%  "
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("E51").Value = "  -3.20%  "
